# Case_0_8 / res_bus / vm_pu.xlsx
# "case with 380 kV done" - slack bus voltage setpoint changed from 1.05 pu to 1.02 pu,
# and the resulting bus voltage magnitudes (columns B-N, rows 2-25) are updated to the
# recomputed power-flow results.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.022927288344651
$ws.Range("D2").Value = 1.025571941809823
$ws.Range("E2").Value = 1.023608011796151
$ws.Range("F2").Value = 1.021383749159382
$ws.Range("I2").Value = 1.029462288934039
$ws.Range("J2").Value = 1.028110375872427
$ws.Range("K2").Value = 1.028397121142402
$ws.Range("L2").Value = 1.026438950580855
$ws.Range("M2").Value = 1.024221239959002
$ws.Range("N2").Value = 1.029570410532077

# row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.024048684826391
$ws.Range("D3").Value = 1.026605898340174
$ws.Range("E3").Value = 1.024563195048109
$ws.Range("F3").Value = 1.023145456043784
$ws.Range("I3").Value = 1.029652169397116
$ws.Range("J3").Value = 1.028869139835098
$ws.Range("K3").Value = 1.029237793935099
$ws.Range("L3").Value = 1.027200644748846
$ws.Range("M3").Value = 1.025786774690352
$ws.Range("N3").Value = 1.030330252026604

# row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.024773206626803
$ws.Range("D4").Value = 1.027274147201752
$ws.Range("E4").Value = 1.025180712360282
$ws.Range("F4").Value = 1.024284069180329
$ws.Range("I4").Value = 1.029772460393206
$ws.Range("J4").Value = 1.029358510568884
$ws.Range("K4").Value = 1.029780380835498
$ws.Range("L4").Value = 1.027692349403139
$ws.Range("M4").Value = 1.026798027858688
$ws.Range("N4").Value = 1.030820317722964

# row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.025077535726505
$ws.Range("D5").Value = 1.027554891827847
$ws.Range("E5").Value = 1.025440186497765
$ws.Range("F5").Value = 1.024762433750256
$ws.Range("I5").Value = 1.029822414703959
$ws.Range("J5").Value = 1.029563860695986
$ws.Range("K5").Value = 1.030008154728576
$ws.Range("L5").Value = 1.027898785306581
$ws.Range("M5").Value = 1.027222747713357
$ws.Range("N5").Value = 1.0310259594708

# row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.025128618760434
$ws.Range("D6").Value = 1.027602019165346
$ws.Range("E6").Value = 1.025483745760531
$ws.Range("F6").Value = 1.024842735465101
$ws.Range("I6").Value = 1.029830766136445
$ws.Range("J6").Value = 1.029598317575852
$ws.Range("K6").Value = 1.030046379690139
$ws.Range("L6").Value = 1.027933430619906
$ws.Range("M6").Value = 1.027294036193906
$ws.Range("N6").Value = 1.031060465283387

# row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.024777274104927
$ws.Range("D7").Value = 1.027277899258343
$ws.Range("E7").Value = 1.025184179975881
$ws.Range("F7").Value = 1.024290462310848
$ws.Range("I7").Value = 1.029773130305399
$ws.Range("J7").Value = 1.029361255963539
$ws.Range("K7").Value = 1.029783425654436
$ws.Range("L7").Value = 1.027695108895282
$ws.Range("M7").Value = 1.026803704590485
$ws.Range("N7").Value = 1.030823067016394

# row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.023306498244405
$ws.Range("D8").Value = 1.025921537060386
$ws.Range("E8").Value = 1.023930935301953
$ws.Range("F8").Value = 1.021979407334314
$ws.Range("I8").Value = 1.029526993140487
$ws.Range("J8").Value = 1.028367136137509
$ws.Range("K8").Value = 1.028681518256834
$ws.Range("L8").Value = 1.026696610314828
$ws.Range("M8").Value = 1.02475068779104
$ws.Range("N8").Value = 1.029827535426183

# row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.020706281754133
$ws.Range("D9").Value = 1.023525312079593
$ws.Range("E9").Value = 1.021718273330448
$ws.Range("F9").Value = 1.017896414808226
$ws.Range("I9").Value = 1.029073537489994
$ws.Range("J9").Value = 1.026603032303448
$ws.Range("K9").Value = 1.026729128211225
$ws.Range("L9").Value = 1.024928154783473
$ws.Range("M9").Value = 1.021119202265194
$ws.Range("N9").Value = 1.028060926362333

# row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.018966900126279
$ws.Range("D10").Value = 1.02192356973283
$ws.Range("E10").Value = 1.02024018068139
$ws.Range("F10").Value = 1.015166663303237
$ws.Range("I10").Value = 1.02875794551203
$ws.Range("J10").Value = 1.025418550448312
$ws.Range("K10").Value = 1.025420227326463
$ws.Range("L10").Value = 1.023743048672178
$ws.Range("M10").Value = 1.018688386508969
$ws.Range("N10").Value = 1.0268747624071

# row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.018212288293181
$ws.Range("D11").Value = 1.021228957279986
$ws.Range("E11").Value = 1.019599418020629
$ws.Range("F11").Value = 1.013982672961656
$ws.Range("I11").Value = 1.028618132993925
$ws.Range("J11").Value = 1.024903634119623
$ws.Range("K11").Value = 1.024851697414014
$ws.Range("L11").Value = 1.023228405458526
$ws.Range("M11").Value = 1.017633364626644
$ws.Range("N11").Value = 1.026359114838163

# row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.017931770716326
$ws.Range("D12").Value = 1.020970787191722
$ws.Range("E12").Value = 1.019361297027874
$ws.Range("F12").Value = 1.013542575291912
$ws.Range("I12").Value = 1.028565725158173
$ws.Range("J12").Value = 1.024712064133843
$ws.Range("K12").Value = 1.024640252040123
$ws.Range("L12").Value = 1.023037018863003
$ws.Range("M12").Value = 1.017241101995887
$ws.Range("N12").Value = 1.026167272801029

# row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.017991952753905
$ws.Range("D13").Value = 1.021026172840917
$ws.Range("E13").Value = 1.019412379948096
$ws.Range("F13").Value = 1.013636992000053
$ws.Range("I13").Value = 1.028576988332716
$ws.Range("J13").Value = 1.024753170484616
$ws.Range("K13").Value = 1.024685619961772
$ws.Range("L13").Value = 1.023078082162989
$ws.Range("M13").Value = 1.017325261051804
$ws.Range("N13").Value = 1.026208437527536

# row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.018189105154877
$ws.Range("D14").Value = 1.021207620141778
$ws.Range("E14").Value = 1.019579737187488
$ws.Range("F14").Value = 1.013946300763705
$ws.Range("I14").Value = 1.028613810644325
$ws.Range("J14").Value = 1.024887805164709
$ws.Range("K14").Value = 1.024834224756481
$ws.Range("L14").Value = 1.023212589989344
$ws.Range("M14").Value = 1.017600947899239
$ws.Range("N14").Value = 1.026343263404317

# row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.018310547844571
$ws.Range("D15").Value = 1.021319394515842
$ws.Range("E15").Value = 1.01968283643779
$ws.Range("F15").Value = 1.014136834538925
$ws.Range("I15").Value = 1.028636435092751
$ws.Range("J15").Value = 1.024970717261792
$ws.Range("K15").Value = 1.024925749478863
$ws.Range("L15").Value = 1.023295434794953
$ws.Range("M15").Value = 1.017770756923871
$ws.Range("N15").Value = 1.026426293246088

# row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.019016950413652
$ws.Range("D16").Value = 1.021969646573149
$ws.Range("E16").Value = 1.020282690235083
$ws.Range("F16").Value = 1.015245197897495
$ws.Range("I16").Value = 1.028767157797095
$ws.Range("J16").Value = 1.025452680814298
$ws.Range("K16").Value = 1.025457921334815
$ws.Range("L16").Value = 1.023777172387241
$ws.Range("M16").Value = 1.018758352053859
$ws.Range("N16").Value = 1.02690894124212

# row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.019459667866326
$ws.Range("D17").Value = 1.022377250163266
$ws.Range("E17").Value = 1.020658763325177
$ws.Range("F17").Value = 1.015939903766933
$ws.Range("I17").Value = 1.028848310480972
$ws.Range("J17").Value = 1.02575445916613
$ws.Range("K17").Value = 1.025791263763576
$ws.Range("L17").Value = 1.024078954700579
$ws.Range("M17").Value = 1.019377177655505
$ws.Range("N17").Value = 1.027211148153827

# row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.019717757989802
$ws.Range("D18").Value = 1.022614897565367
$ws.Range("E18").Value = 1.020878049266079
$ws.Range("F18").Value = 1.016344922762637
$ws.Range("I18").Value = 1.02889534070405
$ws.Range("J18").Value = 1.025930285882797
$ws.Range("K18").Value = 1.025985526429366
$ws.Range("L18").Value = 1.024254836139814
$ws.Range("M18").Value = 1.019737891269928
$ws.Range("N18").Value = 1.027387224564599

# row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.019805736501898
$ws.Range("D19").Value = 1.022695912159333
$ws.Range("E19").Value = 1.020952808106289
$ws.Range("F19").Value = 1.016482991705884
$ws.Range("I19").Value = 1.028911325119623
$ws.Range("J19").Value = 1.025990205164697
$ws.Range("K19").Value = 1.026051736162936
$ws.Range("L19").Value = 1.024314782961638
$ws.Range("M19").Value = 1.019860845427718
$ws.Range("N19").Value = 1.027447228938752

# row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.019412182927333
$ws.Range("D20").Value = 1.02233352861045
$ws.Range("E20").Value = 1.020618421616013
$ws.Range("F20").Value = 1.015865388265888
$ws.Range("I20").Value = 1.028839635094122
$ws.Range("J20").Value = 1.025722101444238
$ws.Range("K20").Value = 1.025755516913462
$ws.Range("L20").Value = 1.024046591130219
$ws.Range("M20").Value = 1.019310808106002
$ws.Range("N20").Value = 1.027178744480258

# row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.01813105482027
$ws.Range("D21").Value = 1.021154192875657
$ws.Range("E21").Value = 1.019530457794886
$ws.Range("F21").Value = 1.013855225722268
$ws.Range("I21").Value = 1.02860298051541
$ws.Range("J21").Value = 1.024848167105535
$ws.Range("K21").Value = 1.024790471768484
$ws.Range("L21").Value = 1.023172987020975
$ws.Range("M21").Value = 1.017519775583605
$ws.Range("N21").Value = 1.026303569054553

# row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.017324277464403
$ws.Range("D22").Value = 1.020411770048415
$ws.Range("E22").Value = 1.01884575563589
$ws.Range("F22").Value = 1.012589551742275
$ws.Range("I22").Value = 1.02845143646571
$ws.Range("J22").Value = 1.024296911665265
$ws.Range("K22").Value = 1.024182157436093
$ws.Range("L22").Value = 1.022622413458014
$ws.Range("M22").Value = 1.016391475438751
$ws.Range("N22").Value = 1.025751530768328

# row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.017752088060076
$ws.Range("D23").Value = 1.020805431286449
$ws.Range("E23").Value = 1.019208792146373
$ws.Range("F23").Value = 1.01326068458421
$ws.Range("I23").Value = 1.028532033707435
$ws.Range("J23").Value = 1.024589311993808
$ws.Range("K23").Value = 1.024504784359275
$ws.Range("L23").Value = 1.022914407256121
$ws.Range("M23").Value = 1.016989821428225
$ws.Range("N23").Value = 1.026044346338876

# row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.019433639758103
$ws.Range("D24").Value = 1.022353284809805
$ws.Range("E24").Value = 1.020636650517636
$ws.Range("F24").Value = 1.01589905920296
$ws.Range("I24").Value = 1.028843556069691
$ws.Range("J24").Value = 1.025736723109866
$ws.Range("K24").Value = 1.025771669902527
$ws.Range("L24").Value = 1.024061215276607
$ws.Range("M24").Value = 1.01934079837957
$ws.Range("N24").Value = 1.027193386910329

# row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.02137952772671
$ws.Range("D25").Value = 1.024145534709155
$ws.Range("E25").Value = 1.022290817826016
$ws.Range("F25").Value = 1.018953287827462
$ws.Range("I25").Value = 1.029193107050463
$ws.Range("J25").Value = 1.027060569482102
$ws.Range("K25").Value = 1.027235146826479
$ws.Range("L25").Value = 1.025386416773207
$ws.Range("N25").Value = 1.028519113296261
